# This script re-applies a row-content permutation to the "Artfynd" sheet.
# A set of data rows were reordered (their full row contents were moved to
# different row numbers) while a handful of other rows stayed fixed.
# We snapshot the full contents (columns A..AY) of every affected row BEFORE
# writing anything, then write each snapshot back into its new row position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$maxCol = 51   # Column AY

function Get-RowSnapshot($ws, $row, $maxCol) {
    $vals = @()
    for ($c = 1; $c -le $maxCol; $c++) {
        $vals += $ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-CellValue($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    if ($val -is [string]) {
        # Force text format so date/time-looking strings (e.g. "2026-02-05")
        # are not auto-converted to Excel date/number serials.
        $cell.NumberFormat = "@"
    }
    $cell.Value2 = $val
}

function Set-RowFromSnapshot($ws, $row, $maxCol, $vals) {
    for ($c = 1; $c -le $maxCol; $c++) {
        Set-CellValue $ws $row $c $vals[$c - 1]
    }
}

# Mapping: new row number -> source (old) row number whose full content
# should be placed there.
$sourceRow = @{}
$sourceRow[2] = 4
$sourceRow[3] = 2
$sourceRow[4] = 3
$sourceRow[7] = 8
$sourceRow[8] = 7
$sourceRow[12] = 13
$sourceRow[13] = 14
$sourceRow[14] = 12
$sourceRow[18] = 19
$sourceRow[19] = 20
$sourceRow[20] = 18

# Snapshot every row that is a source (covers every affected row, since this
# permutation is a union of cycles over the same row set) before any writes.
$snapshots = @{}
foreach ($newRow in $sourceRow.Keys) {
    $oldRow = $sourceRow[$newRow]
    if (-not $snapshots.ContainsKey($oldRow)) {
        $snapshots[$oldRow] = Get-RowSnapshot $ws $oldRow $maxCol
    }
}

# Now write the snapshots into their destination rows.
foreach ($newRow in $sourceRow.Keys) {
    $oldRow = $sourceRow[$newRow]
    Set-RowFromSnapshot $ws $newRow $maxCol $snapshots[$oldRow]
}
